# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de status sheets to reflect the new
# handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 18:53:36"
$wsZhCn.Range("H2").Value = "2016-03-12 18:53:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 18:53:39"
$wsDeDe.Range("H2").Value = "2016-03-12 18:53:58"
